# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (H) and "Correspond Handback
# DateTime" (K) timestamps on the first data row (124a0435-...-md) of the
# per-locale handback sheets, reflecting a newly generated handback report.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("H2").Value = "2016-09-02 10:57:03"
$ws_zhcn.Range("K2").Value = "2016-09-02 10:57:31"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("H2").Value = "2016-09-02 10:57:13"
$ws_dede.Range("K2").Value = "2016-09-02 10:57:38"
